$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new "Tutorial - Gate Optimization and Buffer Design" entry ---
# C10 picks up the same left-aligned body style ("s=5") used throughout column C;
# copy that formatting from C9 before writing the new value.
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C10").Value = "Tutorial - Gate Optimization and Buffer Design"

# New Zoom recording link for the 2/6 lecture tutorial
$ws.Hyperlinks.Add($ws.Range("H10"), "https://iu.zoom.us/rec/share/TAu5wQAkgguM92x84qoZu5pTdF9lLYijKJ2AXb-5RaLxx2F2u3TnhbwLHJfoRU2F.3L3P76JRY2GSMT1B?startTime=1707243575000") | Out-Null
$ws.Range("H9").Copy()
$ws.Range("H10").PasteSpecial(-4122)  # xlPasteFormats -> matches hyperlink style (s=8) used by H8/H9

# --- Row 12: Quiz 1 folded into the HW 1 due date ---
$ws.Range("G12").Value = "HW 1/Quiz 1"

# --- Row 5: Tutorial/NAND2 row now shows "NA" under Slides (D5) ---
$ws.Range("D5").Value = "NA"

# --- Row 10: Slides column (D10) also shows "NA" ---
$ws.Range("C9").Copy()
$ws.Range("D10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D10").Value = "NA"

# --- Row 9: "5.0 - DC and Transient Analysis" lecture; add Slides value 5; drop "HW 1" from E9 ---
$ws.Range("C9").Copy()
$ws.Range("D9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D9").Value = 5
$ws.Range("E9").ClearContents()

# HW 1 assignment moves down into row 10 (was previously shown against row 9)
$ws.Range("E10").Value = "HW 1"

# --- Row 11: Lab 1 due date moves from column G to column E ---
$ws.Range("E11").Value = "Lab 1"
$ws.Range("G11").ClearContents()

# --- widen column C so the longer topic text continues to fit ---
$ws.Columns("C").ColumnWidth = 37.83

$excel.CutCopyMode = $false
